# This script re-applies the row-content permutation described by the commit:
# the identifying "id" (and all associated odds data in columns B:AC) for a
# number of match rows got corrected/reshuffled, while the row's own serial
# number in column A (and its row position) stays put.
#
# Mapping: destRow -> srcRow  (the data that should end up in destRow is the
# data that currently/originally lives in srcRow).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    50 = 51; 51 = 56; 52 = 50; 53 = 55; 54 = 53; 55 = 52; 56 = 54;
    87 = 88; 88 = 87;
    93 = 94; 94 = 93;
    107 = 108; 108 = 107; 109 = 110; 110 = 109;
    123 = 124; 124 = 123;
    135 = 137; 136 = 135; 137 = 136;
    139 = 141; 141 = 139;
    152 = 153; 153 = 152;
    155 = 156; 156 = 155;
    174 = 175; 175 = 174;
    203 = 204; 204 = 203;
    206 = 207; 207 = 206;
    214 = 215; 215 = 214;
    216 = 217; 217 = 216;
    228 = 230; 230 = 228;
}

$firstCol = 2   # column B
$lastCol  = 29  # column AC

# Collect the set of distinct rows involved (sources and destinations are the
# same set here, since every group is a permutation of itself).
$rows = New-Object System.Collections.Generic.HashSet[int]
foreach ($k in $mapping.Keys) { [void]$rows.Add([int]$k) }
foreach ($v in $mapping.Values) { [void]$rows.Add([int]$v) }

# Snapshot the current (pre-edit) B:AC values of every involved row first, so
# that rows which depend on each other cyclically (e.g. 50 <- 51 <- 56 <- 54
# <- 53 <- 55 <- 52 <- 50) are all computed from the ORIGINAL data, not data
# that has already been partially overwritten.
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write each destination row's B:AC range using the snapshot of its
# mapped source row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
